# This script applies the edit described by the diff:
# - Two brand-new data rows are inserted right before the existing row 497
#   (pushing the former rows 497..579 down to 499..581).
# - The two new rows contain fresh "Apio" (celery) price records.
# The sheet's dimension (A1:R579 -> A1:R581) is updated automatically by Excel
# when rows are inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 497 (each Insert() shifts existing row 497
# and everything below it down by one row).
$ws.Rows.Item(497).Insert()
$ws.Rows.Item(497).Insert()

# ---- Fill in the first new row (497) ----
$ws.Range("A497").Value = 10
$ws.Range("B497").Value = "Vega Modelo de Temuco"
$ws.Range("C497").Value = "La Araucanía"
$ws.Range("D497").Value = 45218
$ws.Range("E497").Value = 9
$ws.Range("F497").Value = 100112017
$ws.Range("G497").Value = "Apio"
$ws.Range("H497").Value = "Americana (o)"
$ws.Range("I497").Value = "Primera"
$ws.Range("J497").Value = 110
$ws.Range("K497").Value = 8000
$ws.Range("L497").Value = 8000
$ws.Range("M497").Value = 8000
$ws.Range("N497").Value = "`$/caja 8 unidades"
$ws.Range("O497").Value = "Provincia del Elquí"
$ws.Range("P497").Value = 8000
$ws.Range("Q497").Value = 1
$ws.Range("R497").Value = "Hortaliza"

# ---- Fill in the second new row (498) ----
$ws.Range("A498").Value = 10
$ws.Range("B498").Value = "Vega Modelo de Temuco"
$ws.Range("C498").Value = "La Araucanía"
$ws.Range("D498").Value = 45218
$ws.Range("E498").Value = 9
$ws.Range("F498").Value = 100112017
$ws.Range("G498").Value = "Apio"
$ws.Range("H498").Value = "Americana (o)"
$ws.Range("I498").Value = "Primera"
$ws.Range("J498").Value = 190
$ws.Range("K498").Value = 8000
$ws.Range("L498").Value = 10000
$ws.Range("M498").Value = 9316
$ws.Range("N498").Value = "`$/docena de matas"
$ws.Range("O498").Value = "Provincia del Elquí"
$ws.Range("P498").Value = 1553
$ws.Range("Q498").Value = 6
$ws.Range("R498").Value = "Hortaliza"
